$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns are treated as text so values like "1.001" or
# "20.002.87" are not auto-converted into numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '20.002.87'
$ws.Range("E2").Value = '  -7.06%  '

$ws.Range("D3").Value = '1.411.32'
$ws.Range("E3").Value = '  -7.77%  '

$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").Value = '1.002'
$ws.Range("E5").Value = '  +0.04%  '

$ws.Range("D6").Value = '275.21'
$ws.Range("E6").Value = '  -4.58%  '

$ws.Range("D7").Value = '0.3671'
$ws.Range("E7").Value = '  -5.30%  '

$ws.Range("D8").Value = '0.3118'
$ws.Range("E8").Value = '  -1.51%  '

$ws.Range("D9").Value = '39.81'
$ws.Range("E9").Value = '  -6.73%  '

$ws.Range("D10").Value = '1.033'
$ws.Range("E10").Value = '  -3.09%  '

$ws.Range("D11").Value = '0.06495'
$ws.Range("E11").Value = '  -9.11%  '

$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  -0.01%  '

$ws.Range("D13").Value = '5.480'
$ws.Range("E13").Value = '  -4.19%  '

$ws.Range("D14").Value = '17.70'
$ws.Range("E14").Value = '  -2.20%  '

$ws.Range("D15").Value = '6.182'
$ws.Range("E15").Value = '  -5.55%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '0.00001020'
$ws.Range("E16").Value = '  -6.10%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '1.412.60'
$ws.Range("E17").Value = '  -8.05%  '

$ws.Range("D18").Value = '0.05693'
$ws.Range("E18").Value = '  -13.98%  '

$ws.Range("E19").Value = '  +0.10%  '

$ws.Range("D20").Value = '70.91'
$ws.Range("E20").Value = '  -15.02%  '

$ws.Range("D21").Value = '5.628'
$ws.Range("E21").Value = '  -7.57%  '

$ws.Range("D22").Value = '14.69'
$ws.Range("E22").Value = '  -4.44%  '

$ws.Range("D23").Value = '11.03'
$ws.Range("E23").Value = '  +2.24%  '

$ws.Range("E24").Value = '  -4.74%  '

$ws.Range("D25").Value = '19.997.98'
$ws.Range("E25").Value = '  -7.07%  '

$ws.Range("D26").Value = '2.263'
$ws.Range("E26").Value = '  -4.35%  '

$ws.Range("D27").Value = '133.46'
$ws.Range("E27").Value = '  -10.60%  '

$ws.Range("D28").Value = '17.07'
$ws.Range("E28").Value = '  -6.87%  '

$ws.Range("D29").Value = '1.569.93'
$ws.Range("E29").Value = '  -7.99%  '

$ws.Range("D30").Value = '109.54'
$ws.Range("E30").Value = '  -5.97%  '

$ws.Range("D31").Value = '3.963'
$ws.Range("E31").Value = '  -18.05%  '

$ws.Range("D32").Value = '5.297'
$ws.Range("E32").Value = '  -12.10%  '

$ws.Range("D33").Value = '0.8215'
$ws.Range("E33").Value = '  -13.30%  '

$ws.Range("D34").Value = '0.07689'
$ws.Range("E34").Value = '  -3.70%  '

$ws.Range("D35").Value = '8.428'
$ws.Range("E35").Value = '  -1.05%  '

$ws.Range("D36").Value = '1.483'
$ws.Range("E36").Value = '  -0.17%  '

$ws.Range("D37").Value = '0.05894'
$ws.Range("E37").Value = '  +0.30%  '

$ws.Range("D38").Value = '4.893'
$ws.Range("E38").Value = '  -4.98%  '

$ws.Range("E39").Value = '  +0.11%  '

$ws.Range("D40").Value = '0.02073'
$ws.Range("E40").Value = '  -5.80%  '

$ws.Range("D41").Value = '10.52'
$ws.Range("E41").Value = '  -6.65%  '

$ws.Range("D42").Value = '0.1904'
$ws.Range("E42").Value = '  -5.68%  '

$ws.Range("E43").Value = '  -7.72%  '

$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '0.5302'
$ws.Range("E44").Value = '  -7.53%  '

$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").Value = '3.531'
$ws.Range("E45").Value = '  -4.87%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '12.25'
$ws.Range("E46").Value = '  -7.31%  '

$ws.Range("D47").Value = '0.5191'
$ws.Range("E47").Value = '  -6.28%  '

$ws.Range("D48").Value = '116.07'
$ws.Range("E48").Value = '  +0.56%  '

$ws.Range("D49").Value = '1.768'
$ws.Range("E49").Value = '  -6.26%  '

$ws.Range("D50").Value = '1.040'
$ws.Range("E50").Value = '  -9.97%  '

$ws.Range("D51").Value = '1.002'
$ws.Range("E51").Value = '  +0.00%  '

# Restore the default (unstyled) look for the price/volume columns now
# that the text values have been written, matching the original styling.
$ws.Range("D2:E51").Style = "Normal"
